$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the bottom two rows (old rows 4 and 5)
$ws.Rows("4:5").Delete()

# Update the remaining commessa values to the new ones from the diff
$ws.Range("A2").Value = 252417
$ws.Range("A3").Value = 252980
